$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("All Published Values")
$wsSummary = $wb.Worksheets.Item("Daily Summary")

# Append the new row (row 13) of published values to "All Published Values".
# All source cells are plain text (mirrors the scraper's inlineStr export), so
# values that look like dates/numbers are written with a leading apostrophe to
# force text storage, then the quote-prefix style is reset back to the
# sheet's default (style 0) to match the rest of the table's unstyled cells.
$newRow = 13
$textCols = @{
  1  = "2026-01-02"
  2  = "2026-01-02 20:58:09"
  3  = "697.85"
  4  = "697.85"
  5  = "700.79"
  6  = "700.79"
  7  = "702.88"
  8  = "2026/01/02 20:58:09"
  9  = "2026-01-02 13:09:13"
  10 = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"
}
foreach ($col in $textCols.Keys) {
  $cell = $wsData.Cells.Item($newRow, $col)
  $cell.Value = "'" + $textCols[$col]
  $cell.Style = "Normal"
}

# Extend the autofilter range to cover the newly added row. Re-apply it from
# scratch (turn off, then re-filter the full A1:J13 range) since reassigning
# AutoFilter.Range directly is a no-op here.
$wsData.AutoFilterMode = $false
$null = $wsData.Range("A1:J13").AutoFilter()

# Update the hidden _FilterDatabase defined name for this sheet to the new range
$wb.Names.Item(1).RefersTo = "='All Published Values'!`$A`$1:`$J`$13"

# Update the publishes count on the Daily Summary sheet (11 -> 12)
$wsSummary.Cells.Item(4, 2).Value = 12
